$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.088.96"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "1.888.93"
$ws.Range("E3").Value = "  +1.87%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "307.54"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.22%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.06%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5154"
$c.ClearFormats()
$ws.Range("E7").Value = "  +2.19%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3721"
$c.ClearFormats()
$ws.Range("E8").Value = "  +2.02%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07209"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.80%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9034"
$c.ClearFormats()
$ws.Range("E10").Value = "  +1.53%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "21.00"
$c.ClearFormats()
$ws.Range("E11").Value = "  +1.88%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07616"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("D13").Value = "1.899.61"
$ws.Range("E13").Value = "  +2.39%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "94.80"
$c.ClearFormats()
$ws.Range("E14").Value = "  +3.14%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.271"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.98%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.12%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008497"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.07%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.36"
$c.ClearFormats()
$ws.Range("E18").Value = "  +2.45%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "27.124.88"
$ws.Range("E20").Value = "  +1.02%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.055"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "2.114.82"
$ws.Range("E22").Value = "  +1.09%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.56"
$c.ClearFormats()
$ws.Range("E23").Value = "  +2.45%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.425"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.10%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "145.69"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("E26").Value = "  +0.21%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.ClearFormats()
$ws.Range("E27").Value = "  +1.43%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.161"
$c.ClearFormats()
$ws.Range("E28").Value = "  +5.03%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "114.61"
$c.ClearFormats()
$ws.Range("E29").Value = "  +1.66%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.990"
$c.ClearFormats()
$ws.Range("E30").Value = "  +7.36%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.818"
$c.ClearFormats()
$ws.Range("E31").Value = "  +4.24%  "

$ws.Range("E32").Value = "  +0.19%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05069"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.07%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.197"
$c.ClearFormats()
$ws.Range("E34").Value = "  +4.88%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7584"
$c.ClearFormats()
$ws.Range("E35").Value = "  +2.13%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.015"
$c.ClearFormats()
$ws.Range("E36").Value = "  +0.75%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.274"
$c.ClearFormats()
$ws.Range("E37").Value = "  +1.42%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.563"
$c.ClearFormats()
$ws.Range("E38").Value = "  +2.38%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5623"
$c.ClearFormats()
$ws.Range("E39").Value = "  +5.91%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01992"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.56%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.075"
$c.ClearFormats()
$ws.Range("E41").Value = "  -0.21%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.982"
$c.ClearFormats()
$ws.Range("E42").Value = "  +7.63%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.573"
$c.ClearFormats()
$ws.Range("E43").Value = "  +1.73%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "117.91"
$c.ClearFormats()
$ws.Range("E44").Value = "  -1.27%  "

$ws.Range("E45").Value = "  +3.78%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4794"
$c.ClearFormats()
$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("E47").Value = "  -0.04%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "10.14"
$c.ClearFormats()
$ws.Range("E48").Value = "  +2.18%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.573"
$c.ClearFormats()
$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("E50").Value = "  +0.83%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "63.52"
$c.ClearFormats()
$ws.Range("E51").Value = "  +1.22%  "
